$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "E:\storage\75G.mp4"
$ws.Range("C2").Value = "áda"
$ws.Range("D2").Value = "ádas"
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "E:/New folder\75G.mp4"
